# prod smoke test debugging
# Updates RC-series MSRP rows to MY2021 pricing, renames "RC F TRACK" to
# "RC F FUJI SPEEDWAY EDITION" with refreshed pricing, and appends four new
# "Black Line" special-edition trims (rows 95-98).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$currencyFmt0 = """$""#,##0_);[Red](""$""#,##0)"
$currencyFmt2 = """$""#,##0.00_);[Red](""$""#,##0.00)"

# --- Rows 2-9: RC 300 / RC 350 RWD & AWD trims -> MY2021, refreshed MSRP ---
$ws.Range("C2").Value = 2021
$ws.Range("D2").Value = 42120

$ws.Range("C3").Value = 2021
$ws.Range("D3").Value = 46590

$ws.Range("C4").Value = 2021
$ws.Range("D4").Value = 44810

$ws.Range("C5").Value = 2021
$ws.Range("D5").Value = 48765

$ws.Range("C6").Value = 2021
$ws.Range("D6").Value = 45050

$ws.Range("C7").Value = 2021
$ws.Range("D7").Value = 49520

$ws.Range("C8").Value = 2021
$ws.Range("D8").Value = 47215

$ws.Range("C9").Value = 2021
$ws.Range("D9").Value = 51130

# --- Row 53: RC F -> MY2021, refreshed MSRP ---
$ws.Range("C53").Value = 2021
$ws.Range("D53").Value = 65875

# --- Row 54: RC F TRACK renamed to RC F FUJI SPEEDWAY EDITION, MY2021, refreshed MSRP ---
$ws.Range("B54").Value = "RC F FUJI SPEEDWAY EDITION"
$ws.Range("C54").Value = 2021
$ws.Range("D54").Value = 97100

# --- New rows 95-98: RC 300/350 (AWD) F SPORT Black Line special editions ---
$ws.Range("A95").Value = "9203SE"
$ws.Range("B95").Value = "RC 300 F SPORT Black Line"
$ws.Range("C95").Value = 2021
$ws.Range("D95").Value = 49160
$ws.Range("D95").NumberFormat = $currencyFmt0
$ws.Range("E95").Value = 1025
$ws.Range("E95").NumberFormat = $currencyFmt2

$ws.Range("A96").Value = "9207SE"
$ws.Range("B96").Value = "RC 300 AWD F SPORT Black Line"
$ws.Range("C96").Value = 2021
$ws.Range("D96").Value = 51335
$ws.Range("D96").NumberFormat = $currencyFmt0
$ws.Range("E96").Value = 1025
$ws.Range("E96").NumberFormat = $currencyFmt2

$ws.Range("A97").Value = "9213SE"
$ws.Range("B97").Value = "RC 350 F SPORT Black Line"
$ws.Range("C97").Value = 2021
$ws.Range("D97").Value = 52090
$ws.Range("D97").NumberFormat = $currencyFmt0
$ws.Range("E97").Value = 1025
$ws.Range("E97").NumberFormat = $currencyFmt2

$ws.Range("A98").Value = "9217SE"
$ws.Range("B98").Value = "RC 350 AWD F SPORT Black Line"
$ws.Range("C98").Value = 2021
$ws.Range("D98").Value = 53700
$ws.Range("D98").NumberFormat = $currencyFmt0
$ws.Range("E98").Value = 1025
$ws.Range("E98").NumberFormat = $currencyFmt2

# --- View state: scroll/selection moved toward the bottom of the refreshed table ---
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D99").Select()
